# Update Agrp-Mc3r NATMI output with newly computed TPM-based values.
# - Rows for "ECs" and "MuSCs" sending clusters are removed.
# - Remaining rows (FAPs, Inflammatory-Mac, Resolving-Mac) get recalculated
#   expression / specificity values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: FAPs -> Agrp -> Mc3r -> FAPs -----------------------------------
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Agrp"
$ws.Range("C2").Value = "Mc3r"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.8031959999999999
$ws.Range("H2").Value = 2.409588
$ws.Range("I2").Value = 0.3021857029182209
$ws.Range("J2").Value = 0.3021857029182209
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.018433
$ws.Range("N2").Value = 0.055299
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.014805311868
$ws.Range("R2").Value = 0.133247806812
$ws.Range("S2").Value = 0.3021857029182209
$ws.Range("T2").Value = 0.3021857029182209

# --- Row 3: Inflammatory-Mac -> Agrp -> Mc3r -> FAPs -----------------------
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("B3").Value = "Agrp"
$ws.Range("C3").Value = "Mc3r"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.460649333333333
$ws.Range("H3").Value = 4.381948
$ws.Range("I3").Value = 0.5495387744838921
$ws.Range("J3").Value = 0.5495387744838921
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.018433
$ws.Range("N3").Value = 0.055299
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.02692414916133334
$ws.Range("R3").Value = 0.242317342452
$ws.Range("S3").Value = 0.5495387744838921
$ws.Range("T3").Value = 0.5495387744838921

# --- Row 4: Resolving-Mac -> Agrp -> Mc3r -> FAPs --------------------------
$ws.Range("A4").Value = "Resolving-Mac"
$ws.Range("B4").Value = "Agrp"
$ws.Range("C4").Value = "Mc3r"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.3941096666666666
$ws.Range("H4").Value = 1.182329
$ws.Range("I4").Value = 0.148275522597887
$ws.Range("J4").Value = 0.148275522597887
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.018433
$ws.Range("N4").Value = 0.055299
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.007264623485666667
$ws.Range("R4").Value = 0.065381611371
$ws.Range("S4").Value = 0.148275522597887
$ws.Range("T4").Value = 0.148275522597887

# --- Remove the now-obsolete rows 5 and 6 (ECs / MuSCs) --------------------
$ws.Range("A5:T6").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
